$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("H111").Value = 71432020
$ws.Range("I111").Value = 2109.5715
$ws.Range("K111").Value = 6328.7145
$ws.Range("M111").Value = -3261.7145
$ws.Range("H113").Value = 17513.25
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 19600.857
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 19600.857
$ws.Range("M113").Value = 354
$ws.Range("N113").Value = -26108.857
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("N117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("H121").Value = 1919.2667
$ws.Range("I121").Value = 797.5
$ws.Range("J121").Value = 2327.182
$ws.Range("K121").Value = 2392.5
$ws.Range("L121").Value = 6981.545999999999
$ws.Range("M121").Value = -645.5
$ws.Range("N121").Value = -10475.546
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("H126").Value = 15719.474
$ws.Range("J126").Value = 15719.474
$ws.Range("L126").Value = 15719.474
$ws.Range("N126").Value = -25599.474
$ws.Range("H128").Value = 18810.625
$ws.Range("J128").Value = 18810.625
$ws.Range("L128").Value = 18810.625
$ws.Range("N128").Value = -28770.625
$ws.Range("H129").Value = 1032.0212
$ws.Range("I129").Value = 483.33334
$ws.Range("J129").Value = 1112.317
$ws.Range("K129").Value = 1450.00002
$ws.Range("L129").Value = 3336.951
$ws.Range("M129").Value = 3549.99998
$ws.Range("N129").Value = -13336.951
$ws.Range("H137").Value = 1693
$ws.Range("I137").Value = 1817.3334
$ws.Range("J137").Value = 1320
$ws.Range("K137").Value = 5452.0002
$ws.Range("L137").Value = 3960
$ws.Range("M137").Value = -2902.0002
$ws.Range("N137").Value = -9060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4757.2
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 7595.3335
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 7595.3335
$ws.Range("M25").Value = -98
$ws.Range("N25").Value = -8399.333500000001
$ws.Range("H61").Value = 3396.6667
$ws.Range("I61").Value = 2546
$ws.Range("J61").Value = 4619.5
$ws.Range("K61").Value = 2546
$ws.Range("L61").Value = 4619.5
$ws.Range("M61").Value = -2334
$ws.Range("N61").Value = -5043.5
$ws.Range("H74").Value = 1383.3448
$ws.Range("I74").Value = 1037.8125
$ws.Range("J74").Value = 1808.6154
$ws.Range("K74").Value = 1037.8125
$ws.Range("L74").Value = 1808.6154
$ws.Range("M74").Value = -163.8125
$ws.Range("N74").Value = -3556.6154
$ws.Range("H77").Value = 1383.3448
$ws.Range("I77").Value = 1037.8125
$ws.Range("J77").Value = 1808.6154
$ws.Range("K77").Value = 5189.0625
$ws.Range("L77").Value = 9043.076999999999
$ws.Range("M77").Value = -821.0625
$ws.Range("N77").Value = -17779.077
$ws.Range("H132").Value = 4646.5
$ws.Range("I132").Value = 2879
$ws.Range("J132").Value = 6414
$ws.Range("K132").Value = 8637
$ws.Range("L132").Value = 19242
$ws.Range("M132").Value = -6107
$ws.Range("N132").Value = -24302
$ws.Range("H136").Value = 3396.6667
$ws.Range("I136").Value = 2546
$ws.Range("J136").Value = 4619.5
$ws.Range("K136").Value = 7638
$ws.Range("L136").Value = 13858.5
$ws.Range("M136").Value = -5088
$ws.Range("N136").Value = -18958.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1449.75
$ws.Range("I11").Value = 1533
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 1533
$ws.Range("L11").Value = 1200
$ws.Range("M11").Value = -1393
$ws.Range("N11").Value = -1480
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("H60").Value = 60000
$ws.Range("J60").Value = 60000
$ws.Range("L60").Value = 60000
$ws.Range("N60").Value = -61198
$ws.Range("H107").Value = 834
$ws.Range("I107").Value = 834
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 834
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = 1086
$ws.Range("M107").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H134").Value = 2759.5
$ws.Range("I134").Value = 1564.4615
$ws.Range("J134").Value = 7938
$ws.Range("K134").Value = 4693.3845
$ws.Range("L134").Value = 23814
$ws.Range("M134").Value = -2158.3845
$ws.Range("N134").Value = -28884
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 15645.6
$ws.Range("I35").Value = 1066.6666
$ws.Range("J35").Value = 37514
$ws.Range("K35").Value = 1066.6666
$ws.Range("L35").Value = 37514
$ws.Range("M35").Value = -772.6666
$ws.Range("N35").Value = -38102
$ws.Range("H58").Value = 873.5
$ws.Range("I58").Value = 855.5
$ws.Range("J58").Value = 999.5
$ws.Range("K58").Value = 855.5
$ws.Range("L58").Value = 999.5
$ws.Range("M58").Value = -652.5
$ws.Range("N58").Value = -1405.5
$ws.Range("H93").Value = 17551.75
$ws.Range("I93").Value = 8469
$ws.Range("J93").Value = 44800
$ws.Range("K93").Value = 8469
$ws.Range("L93").Value = 44800
$ws.Range("M93").Value = -6597
$ws.Range("N93").Value = -48544
$ws.Range("H132").Value = 3424.9443
$ws.Range("I132").Value = 3193
$ws.Range("K132").Value = 9579
$ws.Range("M132").Value = -7049
$ws.Range("H134").Value = 3795.7104
$ws.Range("I134").Value = 4871.5654
$ws.Range("J134").Value = 2146.0667
$ws.Range("K134").Value = 14614.6962
$ws.Range("L134").Value = 6438.2001
$ws.Range("M134").Value = -12079.6962
$ws.Range("N134").Value = -11508.2001
$ws.Range("H136").Value = 873.5
$ws.Range("I136").Value = 855.5
$ws.Range("J136").Value = 999.5
$ws.Range("K136").Value = 2566.5
$ws.Range("L136").Value = 2998.5
$ws.Range("M136").Value = -16.5
$ws.Range("N136").Value = -8098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 133.42857
$ws.Range("I8").Value = 133.42857
$ws.Range("K8").Value = 400.28571
$ws.Range("M8").Value = -261.28571
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496
$ws.Range("H131").Value = 1540.6
$ws.Range("I131").Value = 1033.3334
$ws.Range("J131").Value = 1556.2887
$ws.Range("K131").Value = 3100.0002
$ws.Range("L131").Value = 4668.8661
$ws.Range("M131").Value = 1939.9998
$ws.Range("N131").Value = -14748.8661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("H132").Value = 6287.4287
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 7002
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 21006
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -26066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("H32").Value = 6771
$ws.Range("I32").Value = 156.5
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 156.5
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 160.5
$ws.Range("N32").Value = -20634
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2536.3333
$ws.Range("I132").Value = 1575.2858
$ws.Range("J132").Value = 3571.3076
$ws.Range("K132").Value = 4725.857400000001
$ws.Range("L132").Value = 10713.9228
$ws.Range("M132").Value = -2195.857400000001
$ws.Range("N132").Value = -15773.9228
$ws.Range("H136").Value = 3295.673
$ws.Range("I136").Value = 357.66666
$ws.Range("J136").Value = 5285.9355
$ws.Range("K136").Value = 1072.99998
$ws.Range("L136").Value = 15857.8065
$ws.Range("M136").Value = 1477.00002
$ws.Range("N136").Value = -20957.8065
